$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the previously-filled-in employee info (Name / Phone / Mailing / Email) ---
$ws.Range("D4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("J5").Value = ""

# --- Turn the first sample expense line into a generic template example ---
# B13 is a date-formatted cell; typing a date-looking string directly would
# make Excel store it as a real date serial number instead of literal text.
# Stage the literal text on a scratch cell formatted as Text, then copy just
# the value into B13 so its existing number format/style is left untouched.
$scratch = $ws.Range("W54")
$scratch.NumberFormat = "@"
$scratch.Value = "08/08/1996"
$scratch.Copy()
$ws.Range("B13").PasteSpecial(-4163)
$scratch.ClearContents()
$scratch.ClearFormats()

$ws.Range("C13").Value = "Example"

$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = ""
